$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the "Tenure in years" header: shift the remaining header labels
# (Reporting Managers .. HR Comments) one column to the left, then clear
# the now-trailing cell.
$ws.Range("C1").Value = "Reporting Managers"
$ws.Range("D1").Value = "Division"
$ws.Range("E1").Value = "Department"
$ws.Range("F1").Value = "Final Average KRA Grade"
$ws.Range("G1").Value = "Absent Days"
$ws.Range("H1").Value = "HR Comments"
$ws.Range("I1").ClearContents()

# Restore the selection to match the post-edit header range.
$ws.Range("C1:H1").Select()
